$d = $word.ActiveDocument

# Update the date line at the top of the document (wdReplaceOne = 1 to be precise and safe)
$d.Content.Find.Execute("2024-08-29 Thursday", $true, $false, $false, $false, $false, $true, 0, $false, "2024-08-30 Friday", 1) | Out-Null

# Update the answer table (5 data rows x 5 columns).
# Each replacement is scoped to a single cell and uses wdReplaceOne (1) so that
# duplicate text values appearing in different cells (e.g. "10÷5=2, 0") are not
# all replaced by a single Find/Execute call.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$cell.Range.Find.Execute("43÷3=14, 1", $true, $false, $false, $false, $false, $true, 0, $false, "61÷2=30, 1", 1) | Out-Null
$cell = $t.Cell(1, 2)
$cell.Range.Find.Execute("98÷3=32, 2", $true, $false, $false, $false, $false, $true, 0, $false, "84÷3=28, 0", 1) | Out-Null
$cell = $t.Cell(1, 3)
$cell.Range.Find.Execute("97÷2=48, 1", $true, $false, $false, $false, $false, $true, 0, $false, "84÷7=12, 0", 1) | Out-Null
$cell = $t.Cell(1, 4)
$cell.Range.Find.Execute("25÷3=8, 1", $true, $false, $false, $false, $false, $true, 0, $false, "23÷2=11, 1", 1) | Out-Null
$cell = $t.Cell(1, 5)
$cell.Range.Find.Execute("87÷8=10, 7", $true, $false, $false, $false, $false, $true, 0, $false, "97÷7=13, 6", 1) | Out-Null

$cell = $t.Cell(5, 1)
$cell.Range.Find.Execute("51÷3=17, 0", $true, $false, $false, $false, $false, $true, 0, $false, "24÷8=3, 0", 1) | Out-Null
$cell = $t.Cell(5, 2)
$cell.Range.Find.Execute("69÷5=13, 4", $true, $false, $false, $false, $false, $true, 0, $false, "80÷5=16, 0", 1) | Out-Null
$cell = $t.Cell(5, 3)
$cell.Range.Find.Execute("41÷5=8, 1", $true, $false, $false, $false, $false, $true, 0, $false, "65÷9=7, 2", 1) | Out-Null
$cell = $t.Cell(5, 4)
$cell.Range.Find.Execute("31÷2=15, 1", $true, $false, $false, $false, $false, $true, 0, $false, "18÷6=3, 0", 1) | Out-Null
$cell = $t.Cell(5, 5)
$cell.Range.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, 0, $false, "39÷5=7, 4", 1) | Out-Null

$cell = $t.Cell(9, 1)
$cell.Range.Find.Execute("61÷5=12, 1", $true, $false, $false, $false, $false, $true, 0, $false, "93÷6=15, 3", 1) | Out-Null
$cell = $t.Cell(9, 2)
$cell.Range.Find.Execute("87÷6=14, 3", $true, $false, $false, $false, $false, $true, 0, $false, "38÷2=19, 0", 1) | Out-Null
$cell = $t.Cell(9, 3)
$cell.Range.Find.Execute("33÷8=4, 1", $true, $false, $false, $false, $false, $true, 0, $false, "10÷6=1, 4", 1) | Out-Null
$cell = $t.Cell(9, 4)
$cell.Range.Find.Execute("28÷8=3, 4", $true, $false, $false, $false, $false, $true, 0, $false, "20÷8=2, 4", 1) | Out-Null
$cell = $t.Cell(9, 5)
$cell.Range.Find.Execute("10÷5=2, 0", $true, $false, $false, $false, $false, $true, 0, $false, "57÷6=9, 3", 1) | Out-Null

$cell = $t.Cell(13, 1)
$cell.Range.Find.Execute("89÷8=11, 1", $true, $false, $false, $false, $false, $true, 0, $false, "59÷8=7, 3", 1) | Out-Null
$cell = $t.Cell(13, 2)
$cell.Range.Find.Execute("38÷9=4, 2", $true, $false, $false, $false, $false, $true, 0, $false, "59÷6=9, 5", 1) | Out-Null
$cell = $t.Cell(13, 3)
$cell.Range.Find.Execute("54÷7=7, 5", $true, $false, $false, $false, $false, $true, 0, $false, "23÷2=11, 1", 1) | Out-Null
$cell = $t.Cell(13, 4)
$cell.Range.Find.Execute("55÷6=9, 1", $true, $false, $false, $false, $false, $true, 0, $false, "85÷7=12, 1", 1) | Out-Null
$cell = $t.Cell(13, 5)
$cell.Range.Find.Execute("25÷7=3, 4", $true, $false, $false, $false, $false, $true, 0, $false, "26÷4=6, 2", 1) | Out-Null

$cell = $t.Cell(17, 1)
$cell.Range.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 0, $false, "12÷8=1, 4", 1) | Out-Null
$cell = $t.Cell(17, 2)
$cell.Range.Find.Execute("33÷3=11, 0", $true, $false, $false, $false, $false, $true, 0, $false, "26÷6=4, 2", 1) | Out-Null
$cell = $t.Cell(17, 3)
$cell.Range.Find.Execute("11÷3=3, 2", $true, $false, $false, $false, $false, $true, 0, $false, "47÷8=5, 7", 1) | Out-Null
$cell = $t.Cell(17, 4)
$cell.Range.Find.Execute("77÷4=19, 1", $true, $false, $false, $false, $false, $true, 0, $false, "22÷2=11, 0", 1) | Out-Null
$cell = $t.Cell(17, 5)
$cell.Range.Find.Execute("50÷4=12, 2", $true, $false, $false, $false, $false, $true, 0, $false, "28÷5=5, 3", 1) | Out-Null

